# Cucumber Lessons Repeating with Jenkins Examples
# Duplicate the last data row (row 95) into a new row 96, matching the
# existing "Create Country" / "PASSED" / "chrome" / "07.04.23" scenario row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole last row (row 95) into a new row 96, preserving values,
# shared-string usage and formatting exactly (Copy/Paste avoids Excel's
# automatic text->date conversion that a plain .Value assignment would
# trigger for the date-like "07.04.23" string in column D).
$ws.Range("A95:D95").Copy($ws.Range("A96:D96"))
